$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.396.41'
$ws.Range("E2").Value = '  -2.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.643.31'
$ws.Range("E3").Value = '  -3.37%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.86'
$ws.Range("E5").Value = '  -0.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.85'
$ws.Range("E6").Value = '  -1.79%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.543'
$ws.Range("E8").Value = '  -0.82%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.642.62'
$ws.Range("E9").Value = '  -3.33%  '

$ws.Range("E10").Value = '  +0.29%  '

$ws.Range("E11").Value = '  +1.35%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.364'
$ws.Range("E12").Value = '  -0.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.23'
$ws.Range("E13").Value = '  -1.94%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.94'
$ws.Range("E14").Value = '  -2.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.123.42'
$ws.Range("E15").Value = '  -3.46%  '

$ws.Range("E16").Value = '  -2.86%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.371.23'
$ws.Range("E17").Value = '  -2.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.626.64'
$ws.Range("E18").Value = '  -3.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.87'
$ws.Range("E19").Value = '  +0.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.85'
$ws.Range("E20").Value = '  +2.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '364.12'
$ws.Range("E21").Value = '  -2.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.39'
$ws.Range("E22").Value = '  -2.93%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.79'
$ws.Range("E23").Value = '  -3.39%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.14'
$ws.Range("E24").Value = '  +11.84%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.01'
$ws.Range("E25").Value = '  -5.82%  '

$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '70.83'
$ws.Range("E27").Value = '  -4.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.778.10'
$ws.Range("E28").Value = '  -3.43%  '

$ws.Range("E29").Value = '  -3.78%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '554.47'
$ws.Range("E31").Value = '  -6.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.03'
$ws.Range("E32").Value = '  -3.23%  '

$ws.Range("E33").Value = '  -3.94%  '

$ws.Range("E34").Value = '  -1.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.133'
$ws.Range("E35").Value = '  -0.20%  '

$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("E37").Value = '  -4.96%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '157.81'
$ws.Range("E38").Value = '  -2.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.41'
$ws.Range("E39").Value = '  -2.21%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.372'
$ws.Range("E40").Value = '  -2.40%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.27'
$ws.Range("E41").Value = '  -4.38%  '

$ws.Range("E42").Value = '  -4.95%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.93'
$ws.Range("E43").Value = '  -0.36%  '

# Row 44: coin swap/update
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.52'
$ws.Range("E44").Value = '  -5.74%  '

# Row 45: coin swap/update
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.13'
$ws.Range("E46").Value = '  -1.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0303'
$ws.Range("E47").Value = '  -3.66%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.596'
$ws.Range("E48").Value = '  -1.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '154.16'
$ws.Range("E49").Value = '  -1.71%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.88'
$ws.Range("E50").Value = '  -2.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.73'
$ws.Range("E51").Value = '  -3.72%  '
